# Doing Updates for Financials
#
# Refreshes the quarterly/annual figures on the "SNE" financial-statement
# worksheet (Income Statement, Balance Sheet and Cash Flow Statement
# sections) with updated source data. Each entry below is a data row;
# values are written left-to-right starting at column D (the most recent
# period) through column J (the oldest period shown).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SNE")

$updates = @(
    @{ Row = 8; StartCol = 4; Values = @(77237600, 68733400, 73275600, 74271600, 70216100, 61431400, 58698600) }
    @{ Row = 9; StartCol = 4; Values = @(52178000, 47149900, 50682500, 51472100, 49631200, 43989000, 46307400) }
    @{ Row = 10; StartCol = 4; Values = @(25059600, 21583500, 22593200, 22799500, 20584800, 17442400, 12391300) }
    @{ Row = 12; StartCol = 4; Values = @(4145000, 4045000, 4232400, 4197500, 4212900, 4281400) }
    @{ Row = 14; StartCol = 4; Values = @(348900, 1013100) }
    @{ Row = 17; StartCol = 4; Values = @(70594500, 66123500, 70616100, 73651900, 69976600, 59383800, 59306800) }
    @{ Row = 18; StartCol = 4; Values = @(6643100, 2609900, 2659500, 619700, 239500, 2047600, -608200) }
    @{ Row = 20; StartCol = 4; Values = @(-201100, -203800, 321800, -47200, 205300, 381800, 68000) }
    @{ Row = 21; StartCol = 4; Values = @(9713100, 5365900, 6575000, 3781800, 3853900, 5838900, 2352100) }
    @{ Row = 22; StartCol = 4; Values = @(122600, 131500, 228600, 213300, 212100, 241000, 211800) }
    @{ Row = 23; StartCol = 4; Values = @(6319400, 2274600, 2752700, 359200, 232700, 2188400, -752000) }
    @{ Row = 24; StartCol = 4; Values = @(1372000, 1121500, 856900, 802100, 855000, 1269200, 2849800) }
    @{ Row = 26; StartCol = 4; Values = @(4947400, 1153200, 1895800, -443000, -622300, 919200, -3601800) }
    @{ Row = 27; StartCol = 4; Values = @(4436800, 662500, 1336000, -1138900, -1160500, 375500, -4128200) }
    @{ Row = 32; StartCol = 4; Values = @(201100, 203800, -321800, 47200, -205300, -381800, -68000) }
    @{ Row = 33; StartCol = 4; Values = @(4436800, 662500, 1336000, -1138900, -1160500, 375500, -4128200) }
    @{ Row = 35; StartCol = 4; Values = @(4436800, 662500, 1336000, -1138900, -1160500, 375500, -4128200) }
    @{ Row = 41; StartCol = 4; Values = @(14340400, 8679700, 8891900, 8582700, 9460100, 7470300, 8087000) }
    @{ Row = 42; StartCol = 4; Values = @(10636500, 9505000, 8555400, 8468000, 7526400, 6306300, 6155500) }
    @{ Row = 43; StartCol = 4; Values = @(10879500, 10644100, 9579200, 10231900, 9904900, 8358700, 8786500) }
    @{ Row = 44; StartCol = 4; Values = @(6264200, 5793100, 6175600, 6015500, 6634800, 6418900, 6391800) }
    @{ Row = 45; StartCol = 4; Values = @(4671400, 4753800, 4736300, 4650900, 4486000, 4410500, 4524200) }
    @{ Row = 46; StartCol = 4; Values = @(46791900, 39375700, 37938400, 37949000, 38012200, 32964700, 33944900) }
    @{ Row = 47; StartCol = 4; Values = @(97234800, 91410600, 83476100, 77123400, 71587900, 66146800, 57128100) }
    @{ Row = 48; StartCol = 4; Values = @(6684800, 6854100, 7420200, 6683100, 6780100, 7788400, 8416200) }
    @{ Row = 49; StartCol = 4; Values = @(9561200, 10004800, 11047300, 10880700, 12361900, 11826500, 9767300) }
    @{ Row = 52; StartCol = 4; Values = @(12079700, 12006200, 10845500, 10506100, 9874800, 9473600, 10936400) }
    @{ Row = 54; StartCol = 4; Values = @(172353000, 159651000, 150727000, 143142000, 138617000, 128468000, 120193000) }
    @{ Row = 57; StartCol = 4; Values = @(4235700, 4880700, 4980700, 5624800, 6444000, 5171800, 6858500) }
    @{ Row = 58; StartCol = 4; Values = @(6523400, 4683400, 3045900, 2002600, 5818800, 3620200, 3709700) }
    @{ Row = 59; StartCol = 4; Values = @(40050600, 37640400, 35643300, 35272700, 33385000, 31629200, 30382900) }
    @{ Row = 60; StartCol = 4; Values = @(50809700, 47204500, 43670000, 42900100, 43243900, 39008400, 40951000) }
    @{ Row = 61; StartCol = 4; Values = @(5636000, 6160400, 5031700, 6437300, 8286500, 8483400, 6890500) }
    @{ Row = 62; StartCol = 4; Values = @(82853200, 77833300, 73713500, 67284200, 61889700, 56793900, 49659800) }
    @{ Row = 66; StartCol = 4; Values = @(145528000, 137076000, 128459000, 122196000, 118203000, 108650000, 101852000) }
    @{ Row = 72; StartCol = 4; Values = @(13021100, 8898700, 8464400, 7356400, 8500000, 9896800, 9803500) }
    @{ Row = 76; StartCol = 4; Values = @(26825000, 22575100, 22268600, 20946400, 20413600, 19818000, 18341200) }
    @{ Row = 81; StartCol = 4; Values = @(4436800, 662500, 1336000, -1138900, -1160500, 375500, -4128200) }
    @{ Row = 83; StartCol = 4; Values = @(3267500, 2956500, 3589700, 3205800, 3405300, 3405700, 2889100) }
    @{ Row = 89; StartCol = 4; Values = @(11344900, 7315700, 6771800, 6821900, 6003600, 4304500, 4696600) }
    @{ Row = 91; StartCol = 4; Values = @(-2377400, -3014900, -3393700, -1951900, -2562500, -2951500, -3458200) }
    @{ Row = 94; StartCol = 4; Values = @(-7432700, -11335900, -9314800, -5782300, -6422900, -6375700, -7981300) }
    @{ Row = 96; StartCol = 4; Values = @(-257500, -228700, -115300, -119000, -231800, -226500, -226700) }
    @{ Row = 100; StartCol = 4; Values = @(2228000, 4088800, 3436300, -2379300, 1879200, 800300, 2326300) }
    @{ Row = 101; StartCol = 4; Values = @(-479500, -280800, -584100, 462300, 529900, 654200, -125000) }
    @{ Row = 102; StartCol = 4; Values = @(5660700, -212200, 309200, -877400, 1989700, -616700, -1083300) }
)

foreach ($u in $updates) {
    $col = $u.StartCol
    foreach ($v in $u.Values) {
        $ws.Cells.Item($u.Row, $col).Value = $v
        $col = $col + 1
    }
}

$ws.Application.Calculate()
